$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: set Title, Date and Description values ---
$meta = $wb.Worksheets.Item("Metadata")

# Row 5 = Title
$meta.Range("B5").Value = "DMI Nom Distributeur"

# Row 8 = Date (refresh publication date/time)
$meta.Range("B8").Value = "2026-02-25T08:15:31+00:00"

# Row 12 = Description
$meta.Range("B12").Value = "Extension créée dans ce volet pour représenter le nom distributeur."

# --- "Elements" sheet: give the root Extension element a Short/Definition ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the top-level "Extension" element; column L = Short, column M = Definition
$elements.Range("L2").Value = "DMI Nom Distributeur"
$elements.Range("M2").Value = "Extension créée dans ce volet pour représenter le nom distributeur."
